$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All existing hyperlink targets are changing (new batch of tweets), so clear them first
$ws.Hyperlinks.Delete()

# Refresh the data rows (A:E) for rows 4-21 with the new tweet set
$ws.Range("A4").Value2 = "Sep 09 2017"
$ws.Range("B4").Value2 = "negative"
$ws.Range("C4").Value2 = -0.1
$ws.Range("D4").Value2 = "https://twitter.com/statuses/910212924504866817"
$ws.Range("E4").Value2 = "#0daytoday #Tecnovision DLX Spot - Arbitrary File Upload Vulnerability [remote #exploits #Vulnerability #0day… https://t.co/tj8JlLp2NM"

$ws.Range("A5").Value2 = "Sep 09 2017"
$ws.Range("B5").Value2 = "negative"
$ws.Range("C5").Value2 = -0.1
$ws.Range("D5").Value2 = "https://twitter.com/statuses/910212922143522816"
$ws.Range("E5").Value2 = "#0daytoday #Tecnovision DLX Spot - Authentication Bypass Vulnerability [remote #exploits #Vulnerability #0day… https://t.co/wpMvHxrNGV"

$ws.Range("A6").Value2 = "Sep 09 2017"
$ws.Range("B6").Value2 = "negative"
$ws.Range("C6").Value2 = -0.1
$ws.Range("D6").Value2 = "https://twitter.com/statuses/910212919547240448"
$ws.Range("E6").Value2 = "#0daytoday #Tecnovision DLX Spot - SSH Backdoor Vulnerability CVE-2017-12928 [remote #exploits #Vulnerability #0day… https://t.co/k1z5DvTXx3"

$ws.Range("A7").Value2 = "Sep 09 2017"
$ws.Range("B7").Value2 = "negative"
$ws.Range("C7").Value2 = -0.1
$ws.Range("D7").Value2 = "https://twitter.com/statuses/910204857511903232"
$ws.Range("E7").Value2 = "#0daytoday #HPE &amp;lt; 7.2 - Java Deserialization Exploit CVE-2016-4372 [remote #exploits #0day #Exploit] https://t.co/HzChSlUqJA"

$ws.Range("A8").Value2 = "Sep 09 2017"
$ws.Range("B8").Value2 = "negative"
$ws.Range("C8").Value2 = -0.05
$ws.Range("D8").Value2 = "https://twitter.com/statuses/909700296158339072"
$ws.Range("E8").Value2 = "#0daytoday #Netdecision 5.8.2 - Local Privilege Escalation Exploit CVE-2017-14311 [remote #exploits #0day #Exploit] https://t.co/hBu3SRhWx5"

$ws.Range("A9").Value2 = "Sep 09 2017"
$ws.Range("B9").Value2 = "negative"
$ws.Range("C9").Value2 = -0.1
$ws.Range("D9").Value2 = "https://twitter.com/statuses/909526400612188160"
$ws.Range("E9").Value2 = "#0daytoday #EMC Data Protection Advisor Hardcoded Password Vulnerability [remote #exploits #Vulnerability #0day… https://t.co/cKPXhkYoxS"

$ws.Range("A10").Value2 = "Sep 09 2017"
$ws.Range("B10").Value2 = "negative"
$ws.Range("C10").Value2 = -0.1
$ws.Range("D10").Value2 = "https://twitter.com/statuses/908810169148157953"
$ws.Range("E10").Value2 = "#0daytoday #Astaro Security Gateway 7 Remote Code Execution #Exploit https://t.co/GoZqzFZIkl"

$ws.Range("A11").Value2 = "Sep 09 2017"
$ws.Range("B11").Value2 = "negative"
$ws.Range("C11").Value2 = -0.1
$ws.Range("D11").Value2 = "https://twitter.com/statuses/908809924876099587"
$ws.Range("E11").Value2 = "#0daytoday #DLink DIR8xx Remote Root Code Execution #Exploit https://t.co/B98M14ZcxE"

$ws.Range("A12").Value2 = "Sep 09 2017"
$ws.Range("B12").Value2 = "negative"
$ws.Range("C12").Value2 = -0.1
$ws.Range("D12").Value2 = "https://twitter.com/statuses/908649382400073728"
$ws.Range("E12").Value2 = "#0daytoday #VLC Media Player iOS App 2.7.8 File Disclosure Vulnerability [remote #exploits #Vulnerability #0day… https://t.co/et1wT2KvKZ"

$ws.Range("A13").Value2 = "Sep 09 2017"
$ws.Range("B13").Value2 = "negative"
$ws.Range("C13").Value2 = -0.1
$ws.Range("D13").Value2 = "https://twitter.com/statuses/908649380026097666"
$ws.Range("E13").Value2 = "#0daytoday #VIPA Automation WinPLC7 5.0.45.5921 Buffer Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/SBIYWWRYIn"

$ws.Range("A14").Value2 = "Sep 09 2017"
$ws.Range("B14").Value2 = "negative"
$ws.Range("C14").Value2 = -0.1
$ws.Range("D14").Value2 = "https://twitter.com/statuses/908649377253662724"
$ws.Range("E14").Value2 = "#0daytoday #Disk Pulse Server 2.2.34 Buffer Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/BujSgn4Tn3"

$ws.Range("A15").Value2 = "Sep 09 2017"
$ws.Range("B15").Value2 = "negative"
$ws.Range("C15").Value2 = -0.1
$ws.Range("D15").Value2 = "https://twitter.com/statuses/908649374812479488"
$ws.Range("E15").Value2 = "#0daytoday #haneWIN DNS Server 1.5.3 Buffer Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/YAFn3NKaW5"

$ws.Range("A16").Value2 = "Sep 09 2017"
$ws.Range("B16").Value2 = "negative"
$ws.Range("C16").Value2 = -0.1
$ws.Range("D16").Value2 = "https://twitter.com/statuses/908649372778246144"
$ws.Range("E16").Value2 = "#0daytoday #KingScada AlarmServer 3.1.2.13 Buffer Overflow Exploit CVE-2014-0787 [remote #exploits #0day #Exploit] https://t.co/w6UILqQNJo"

$ws.Range("A17").Value2 = "Sep 09 2017"
$ws.Range("B17").Value2 = "negative"
$ws.Range("C17").Value2 = -0.1
$ws.Range("D17").Value2 = "https://twitter.com/statuses/908083799296299008"
$ws.Range("E17").Value2 = "#0daytoday #EMC CMCNE Inmservlets.war FileUploadController 11.2.1 - Remote Code Execution #Exploit https://t.co/qB4nS4Wrxe"

$ws.Range("A18").Value2 = "Sep 09 2017"
$ws.Range("B18").Value2 = "negative"
$ws.Range("C18").Value2 = -0.1
$ws.Range("D18").Value2 = "https://twitter.com/statuses/908083683684470784"
$ws.Range("E18").Value2 = "#0daytoday #EMC CMCNE 11.2.1 - FileUploadController Remote Code Execution #Exploit https://t.co/0Ot4rUrdbL"

$ws.Range("A19").Value2 = "Sep 09 2017"
$ws.Range("B19").Value2 = "negative"
$ws.Range("C19").Value2 = -0.1
$ws.Range("D19").Value2 = "https://twitter.com/statuses/908083588989628416"
$ws.Range("E19").Value2 = "#0daytoday #Dameware Mini Remote Control 4.0 - Username Stack Buffer Overflow #Exploit https://t.co/fN7utYLvTi"

$ws.Range("A20").Value2 = "Sep 09 2017"
$ws.Range("B20").Value2 = "negative"
$ws.Range("C20").Value2 = -0.1
$ws.Range("D20").Value2 = "https://twitter.com/statuses/908083451194167296"
$ws.Range("E20").Value2 = "#0daytoday #Cloudview NMS &lt; 2.00b - Arbitrary File Upload #Exploit https://t.co/ItHdb4XWJD"

$ws.Range("A21").Value2 = "Sep 09 2017"
$ws.Range("B21").Value2 = "negative"
$ws.Range("C21").Value2 = -0.1
$ws.Range("D21").Value2 = "https://twitter.com/statuses/908035000179859456"
$ws.Range("E21").Value2 = "#0daytoday #Microsoft #Windows .NET #Framework - Remote Code Execution #0day #Exploit https://t.co/lSpKziNdZs"

# Re-create hyperlinks for D4:D21, matching the existing visual style (underline, blue font)
$ws.Hyperlinks.Add($ws.Range("D4"), "https://twitter.com/statuses/910212924504866817") | Out-Null
$ws.Range("D4").Font.Underline = 2
$ws.Range("D4").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D5"), "https://twitter.com/statuses/910212922143522816") | Out-Null
$ws.Range("D5").Font.Underline = 2
$ws.Range("D5").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D6"), "https://twitter.com/statuses/910212919547240448") | Out-Null
$ws.Range("D6").Font.Underline = 2
$ws.Range("D6").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D7"), "https://twitter.com/statuses/910204857511903232") | Out-Null
$ws.Range("D7").Font.Underline = 2
$ws.Range("D7").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D8"), "https://twitter.com/statuses/909700296158339072") | Out-Null
$ws.Range("D8").Font.Underline = 2
$ws.Range("D8").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D9"), "https://twitter.com/statuses/909526400612188160") | Out-Null
$ws.Range("D9").Font.Underline = 2
$ws.Range("D9").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D10"), "https://twitter.com/statuses/908810169148157953") | Out-Null
$ws.Range("D10").Font.Underline = 2
$ws.Range("D10").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D11"), "https://twitter.com/statuses/908809924876099587") | Out-Null
$ws.Range("D11").Font.Underline = 2
$ws.Range("D11").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D12"), "https://twitter.com/statuses/908649382400073728") | Out-Null
$ws.Range("D12").Font.Underline = 2
$ws.Range("D12").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D13"), "https://twitter.com/statuses/908649380026097666") | Out-Null
$ws.Range("D13").Font.Underline = 2
$ws.Range("D13").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D14"), "https://twitter.com/statuses/908649377253662724") | Out-Null
$ws.Range("D14").Font.Underline = 2
$ws.Range("D14").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D15"), "https://twitter.com/statuses/908649374812479488") | Out-Null
$ws.Range("D15").Font.Underline = 2
$ws.Range("D15").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D16"), "https://twitter.com/statuses/908649372778246144") | Out-Null
$ws.Range("D16").Font.Underline = 2
$ws.Range("D16").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D17"), "https://twitter.com/statuses/908083799296299008") | Out-Null
$ws.Range("D17").Font.Underline = 2
$ws.Range("D17").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D18"), "https://twitter.com/statuses/908083683684470784") | Out-Null
$ws.Range("D18").Font.Underline = 2
$ws.Range("D18").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D19"), "https://twitter.com/statuses/908083588989628416") | Out-Null
$ws.Range("D19").Font.Underline = 2
$ws.Range("D19").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D20"), "https://twitter.com/statuses/908083451194167296") | Out-Null
$ws.Range("D20").Font.Underline = 2
$ws.Range("D20").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("D21"), "https://twitter.com/statuses/908035000179859456") | Out-Null
$ws.Range("D21").Font.Underline = 2
$ws.Range("D21").Font.Color = 16711680

# Drop the auto-registered built-in "Hyperlink" cell style; the sheet already used plain
# direct formatting (font s=2) for link cells before this edit, so keep that convention
# instead of leaving an unused named style behind.
$wb.Styles.Item("Hyperlink").Delete()
